# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the marksheet summary numbers (rows 10-12) ---
# Row labels (A10/A11/A12) pick up the same "mtitleStyle" formatting that is
# already used by the header row above them (A9).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 10
$ws.Range("D10").Value = 18
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 40
$ws.Range("E12").Value = "40/112"

# --- Remove the third (G/H) answer-key block entirely (rows 15-40) ---
$ws.Range("G15:H40").Clear()

# --- Only the first two quiz questions (rows 16-17) keep their second
#     answer-key block (columns D/E); the rest (rows 18-40) lose it. ---
$ws.Range("D18:E40").Clear()

# Row 16/17 column D now mirrors the "Correct Ans" value and takes on the
# same "correctStyle" formatting used in column E.
$ws.Range("E16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("E17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

# --- Mark the rows where the student's answer (column B) equals the
#     correct answer by writing that answer into column A with the
#     "correctStyle" formatting (same style already used on column B). ---
$ws.Range("B22").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option D"

$ws.Range("B27").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Option A"

$ws.Range("B29").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("B31").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = "Option D"

$ws.Range("B32").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Option C"

$ws.Range("B35").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Option D"

$ws.Range("B37").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Option A"

$ws.Range("B39").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = "Option D"
